$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "MEC-2B-Resistencia mecanica" class from row 3 (quarta, B3 & D3)
# to row 6 and row 7 (quarta, D6 & D7) - full logic part of bimestral classes
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("D6").Value = "MEC-2B-Resistencia mecanica"
$ws.Range("D7").Value = "MEC-2B-Resistencia mecanica"
